# "Generate Report for Handback"
#
# The handoff run finished and the handback report now needs to reflect
# that the localized content has been handed back and is back in sync
# with en-US. This updates the status text, stamps the handback
# date/time, and records the "Latest Target File" / "Latest Handback
# File" links (with their handback timestamps) for each row of the
# zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (this text is shared by the Overview summary columns and the
#    per-language Status column)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 2) Latest Handback DateTime (column H): stamp the real handback time
#    in place of the empty-date placeholder
# ---------------------------------------------------------------------
$zhcn.Range("H2").Value = "2016-03-20 20:13:44"
$zhcn.Range("H3").Value = "2016-03-20 20:13:44"
$dede.Range("H2").Value = "2016-03-20 20:13:50"
$dede.Range("H3").Value = "2016-03-20 20:13:50"

# ---------------------------------------------------------------------
# 3) Latest Target File (F) / Latest Handback File (G): record the
#    handed-back target files, linking back to the same source/target
#    locations already used for the handoff record of each row.
# ---------------------------------------------------------------------

# zh-cn, row 2 (cdc6a913-...)
$zhcn.Range("F2").Value = "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d918e349ffb510b38728ca0871914c799a95602c/e2e/cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md", "", "", "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md") | Out-Null

$zhcn.Range("G2").Value = "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ff27559281a054cb643a9185de21deb2a4d48a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.zh-cn.xlf", "", "", "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.zh-cn.xlf") | Out-Null

# zh-cn, row 3 (f8e2b10f-...)
$zhcn.Range("F3").Value = "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d918e349ffb510b38728ca0871914c799a95602c/e2e/f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md", "", "", "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md") | Out-Null

$zhcn.Range("G3").Value = "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ff27559281a054cb643a9185de21deb2a4d48a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.zh-cn.xlf", "", "", "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.zh-cn.xlf") | Out-Null

# de-de, row 2 (cdc6a913-...)
$dede.Range("F2").Value = "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d918e349ffb510b38728ca0871914c799a95602c/e2e/cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md", "", "", "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.md") | Out-Null

$dede.Range("G2").Value = "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36fd1ad7eb904774290a8c1f28c10822f4951a79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.de-de.xlf", "", "", "cdc6a913-a0c2-449f-b726-b2e6c0f0b717.4a73995afbef76dc2dbd41edc9bceb21f87b2c79.de-de.xlf") | Out-Null

# de-de, row 3 (f8e2b10f-...)
$dede.Range("F3").Value = "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/d918e349ffb510b38728ca0871914c799a95602c/e2e/f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md", "", "", "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.md") | Out-Null

$dede.Range("G3").Value = "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36fd1ad7eb904774290a8c1f28c10822f4951a79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.de-de.xlf", "", "", "f8e2b10f-c05e-4ddc-8a3b-e4d848b56c70.6409e77a049f153c1c4f80e749928e632e63dc87.de-de.xlf") | Out-Null

# ---------------------------------------------------------------------
# 4) Match the existing "hyperlink" look (underlined custom blue) that
#    the other link columns (A, B, D) already use on these rows.
# ---------------------------------------------------------------------
foreach ($cellRef in @("F2", "G2", "F3", "G3")) {
    $zhcn.Range($cellRef).Font.Underline = $true
    $zhcn.Range($cellRef).Font.Color = 15570276
    $dede.Range($cellRef).Font.Underline = $true
    $dede.Range($cellRef).Font.Color = 15570276
}
